$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 18: update title and link
$ws.Range("D18").Value = "KoBERT, KoGPT-2"
$ws.Range("E18").Value = "https://freesearch.pe.kr/archives/5579"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] Correcting Deviations from Normality: A Reformulated Diffusion model for Multi-class Unsupervised AD"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?mod=document&uid=3204"
